# Append two new rows (12 and 13) to the end of the sheet data, duplicating
# the content of rows 7 and 8 (Okanogan / Salmon Creek-Green Lake / Salmon 16-4
# and Okanogan / Aeneas Creek Okanogan / Aeneas 16-1 fish passage barrier rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - duplicate of row 7
$ws.Range("A12").Value = "Okanogan"
$ws.Range("B12").Value = "Salmon Creek-Green Lake"
$ws.Range("C12").Value = "Salmon 16-4"
$ws.Range("D12").Value = "Fish Passage Restoration"
$ws.Range("E12").Value = "Fish Passage Barriers"
$ws.Range("F12").Value = "Barriers_pathway"

# Row 13 - duplicate of row 8
$ws.Range("A13").Value = "Okanogan"
$ws.Range("B13").Value = "Aeneas Creek Okanogan"
$ws.Range("C13").Value = "Aeneas 16-1"
$ws.Range("D13").Value = "Fish Passage Restoration"
$ws.Range("E13").Value = "Fish Passage Barriers"
$ws.Range("F13").Value = "Barriers_pathway"
